# The document starts with a single paragraph containing one run:
#   "Hy creating a new Microsoft word document in repo"
# It must become a paragraph with three runs:
#   1. "Original "                              (no special formatting)
#   2. "alternate"                               (Arial 11.5pt, gray shading F2F2F2)
#   3. " document is replacing to decl.dock"     (Arial 11.5pt, gray shading F2F2F2)
#
# Word's Shading object (Range.Shading / Font.Shading) only ever lands on the
# paragraph's <w:pPr><w:shd/></w:pPr>, so to get run-level <w:rPr><w:shd/></w:rPr>
# shading we build the replacement runs as a WordprocessingML fragment and
# splice it in with Range.InsertXML, which replaces the target range's content
# verbatim (including run-level formatting we specify).

$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(1)
$full = $p.Range

# Exclude the trailing paragraph-mark character from the replace range so the
# paragraph itself (and its identity/properties) is kept intact and no extra
# empty paragraph is introduced.
$target = $d.Range($full.Start, $full.End - 1)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$runProps = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2"/></w:rPr>'

$run1 = '<w:r><w:t xml:space="preserve">Original </w:t></w:r>'
$run2 = '<w:r>' + $runProps + '<w:t>alternate</w:t></w:r>'
$run3 = '<w:r>' + $runProps + '<w:t xml:space="preserve"> document is replacing to decl.dock</w:t></w:r>'

$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' `
  + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
  + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
  + '<pkg:xmlData><w:document ' + $ns + '><w:body><w:p>' + $run1 + $run2 + $run3 + '</w:p></w:body></w:document>' `
  + '</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xmlFrag)
